# Append a new order row (row 8) to the Orders sheet, matching the
# data that was appended in the source commit.
#
# Columns: A=name, B=product, C=date, D=address, E=phone, F=quantity
#
# E (phone) and F (quantity) look like numbers ("09372979927", "2") but
# must stay stored as text (same as every other row in the sheet), so we
# temporarily force a Text number format before assigning those values
# and then restore the cell style back to Normal so no stray formatting
# is left behind on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Swapnil Durafe"
$ws.Range("B8").Value = "Paneer"
$ws.Range("C8").Value = "16/7/2025, 12:19:45 am"
$ws.Range("D8").Value = "D/4, C-204, Laxmi Park-1, Lokmanya Nagar, Thane-west"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "09372979927"
$ws.Range("E8").Style = "Normal"

$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2"
$ws.Range("F8").Style = "Normal"
